# P3 - Requirements Stack: fill in missing "Sprint No." values in column E
# and move the active cell selection, per the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 2
$ws.Range("E8").Value = 2
$ws.Range("E9").Value = 2
$ws.Range("E10").Value = 2
$ws.Range("E14").Value = 2
$ws.Range("E15").Value = 2
$ws.Range("E16").Value = 2
$ws.Range("E17").Value = 1
$ws.Range("E18").Value = 1
$ws.Range("E19").Value = 2
$ws.Range("E22").Value = 2
$ws.Range("E23").Value = 1
$ws.Range("E25").Value = 1
$ws.Range("E33").Value = 2
$ws.Range("E34").Value = 2
$ws.Range("E35").Value = 2

# Update the active cell/selection to match the saved view state
$ws.Range("H5").Select()
